$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("D7").Value = 0.03474168402683874
$ws.Range("E7").Value = 0.1458085355288644
$ws.Range("F7").Value = 0.200283817889967
$ws.Range("G7").Value = 0.1863912123111998
$ws.Range("H7").Value = 22.75886895093276
$ws.Range("L7").Value = 0.03271454534848238
$ws.Range("M7").Value = 0.1397425365277118
$ws.Range("N7").Value = 0.1906060551134697
$ws.Range("O7").Value = 0.1808716267093387
$ws.Range("P7").Value = 21.94211006533216

$ws = $wb.Worksheets.Item(2)
$ws.Range("D7").Value = 0.01986410401870377
$ws.Range("E7").Value = 0.1074974094720334
$ws.Range("F7").Value = 0.2286616587492813
$ws.Range("G7").Value = 0.1409400724375568
$ws.Range("H7").Value = 19.04382287636304
$ws.Range("L7").Value = 0.03010178652857923
$ws.Range("M7").Value = 0.1366043386343726
$ws.Range("N7").Value = 0.2169158077340615
$ws.Range("O7").Value = 0.1734986643423494
$ws.Range("P7").Value = 22.58008899310908

$ws = $wb.Worksheets.Item(3)
$ws.Range("D7").Value = 0.02469764282300885
$ws.Range("E7").Value = 0.1371902006883774
$ws.Range("F7").Value = 0.2110627121139214
$ws.Range("G7").Value = 0.1571548370970771
$ws.Range("H7").Value = 23.85666160286808
$ws.Range("L7").Value = 0.02042804453966603
$ws.Range("M7").Value = 0.1141099412225236
$ws.Range("N7").Value = 0.1634396800258065
$ws.Range("O7").Value = 0.1429267103786623
$ws.Range("P7").Value = 17.34919584225236

$ws = $wb.Worksheets.Item(4)
$ws.Range("D7").Value = 0.2268196840211071
$ws.Range("E7").Value = 0.3799973223033592
$ws.Range("F7").Value = 0.1337282000628839
$ws.Range("G7").Value = 0.4762559018228615
$ws.Range("H7").Value = 14.8883871491639
$ws.Range("L7").Value = 0.1406106730781708
$ws.Range("M7").Value = 0.3288962682014027
$ws.Range("N7").Value = 0.1201710028135958
$ws.Range("O7").Value = 0.3749808969509924
$ws.Range("P7").Value = 12.56729371014039

$ws = $wb.Worksheets.Item(5)
$ws.Range("D7").Value = 0.1453884869818583
$ws.Range("E7").Value = 0.2648971741297862
$ws.Range("F7").Value = 0.137453697281873
$ws.Range("G7").Value = 0.3812984224749144
$ws.Range("H7").Value = 11.85807733081818
$ws.Range("L7").Value = 0.1687464144037276
$ws.Range("M7").Value = 0.3071946578586759
$ws.Range("N7").Value = 0.1375970219320574
$ws.Range("O7").Value = 0.4107875538568903
$ws.Range("P7").Value = 13.16750320109533

$ws = $wb.Worksheets.Item(6)
$ws.Range("D7").Value = 0.255749789652545
$ws.Range("E7").Value = 0.3981418332102931
$ws.Range("F7").Value = 0.1628789353654156
$ws.Range("G7").Value = 0.5057171043701656
$ws.Range("H7").Value = 16.25579552736434
$ws.Range("L7").Value = 0.1156087895734882
$ws.Range("M7").Value = 0.2453968391016066
$ws.Range("N7").Value = 0.09936917755563611
$ws.Range("O7").Value = 0.3400129255976722
$ws.Range("P7").Value = 9.935974330173618

$ws = $wb.Worksheets.Item(7)
$ws.Range("D7").Value = 0.02135653082893724
$ws.Range("E7").Value = 0.1153650512726327
$ws.Range("F7").Value = 0.1625405171978774
$ws.Range("G7").Value = 0.1461387382898089
$ws.Range("H7").Value = 18.31261848374522
$ws.Range("L7").Value = 0.02815216176807918
$ws.Range("M7").Value = 0.1338620782158162
$ws.Range("N7").Value = 0.186685914320878
$ws.Range("O7").Value = 0.167786059516514
$ws.Range("P7").Value = 21.52357791571867

$ws = $wb.Worksheets.Item(8)
$ws.Range("D7").Value = 0.03011623746378193
$ws.Range("E7").Value = 0.1201631848612445
$ws.Range("F7").Value = 0.1912492681872061
$ws.Range("G7").Value = 0.1735403050123571
$ws.Range("H7").Value = 19.6840294685524
$ws.Range("L7").Value = 0.02888074148548015
$ws.Range("M7").Value = 0.1172071508330191
$ws.Range("N7").Value = 0.1971851472833981
$ws.Range("O7").Value = 0.1699433478706364
$ws.Range("P7").Value = 19.62628280541363

$ws = $wb.Worksheets.Item(9)
$ws.Range("D7").Value = 0.01669457329228728
$ws.Range("E7").Value = 0.1013542524358794
$ws.Range("F7").Value = 0.2021353547927391
$ws.Range("G7").Value = 0.1292074815646806
$ws.Range("H7").Value = 17.39194282618048
$ws.Range("L7").Value = 0.01688258366375536
$ws.Range("M7").Value = 0.1031553429136347
$ws.Range("N7").Value = 0.1886383453068354
$ws.Range("O7").Value = 0.1299329968243454
$ws.Range("P7").Value = 17.30205941511477

$ws = $wb.Worksheets.Item(10)
$ws.Range("D7").Value = 0.4022691100356651
$ws.Range("E7").Value = 0.5196450310785209
$ws.Range("F7").Value = 0.2153530762739931
$ws.Range("G7").Value = 0.634246884135559
$ws.Range("H7").Value = 22.67085693534974
$ws.Range("L7").Value = 0.4048508920880359
$ws.Range("M7").Value = 0.5409121451494966
$ws.Range("N7").Value = 0.2220282511307675
$ws.Range("O7").Value = 0.6362789420435316
$ws.Range("P7").Value = 23.047178221503

$ws = $wb.Worksheets.Item(11)
$ws.Range("D7").Value = 0.2072395632587745
$ws.Range("E7").Value = 0.3784818606893113
$ws.Range("F7").Value = 0.1517637268855607
$ws.Range("G7").Value = 0.4552357227401805
$ws.Range("H7").Value = 15.08821711678507
$ws.Range("L7").Value = 0.1807129184466426
$ws.Range("M7").Value = 0.3450079260695241
$ws.Range("N7").Value = 0.1407867435407954
$ws.Range("O7").Value = 0.4251034208832511
$ws.Range("P7").Value = 13.78540033198624

$ws = $wb.Worksheets.Item(12)
$ws.Range("D7").Value = 0.1891865924541661
$ws.Range("E7").Value = 0.3684293981733054
$ws.Range("F7").Value = 0.1376300008617787
$ws.Range("G7").Value = 0.4349558511552248
$ws.Range("H7").Value = 15.26015318292194
$ws.Range("L7").Value = 0.159079587431299
$ws.Range("M7").Value = 0.329277868006125
$ws.Range("N7").Value = 0.1238098121261744
$ws.Range("O7").Value = 0.3988478249048114
$ws.Range("P7").Value = 13.72474121846528

$ws = $wb.Worksheets.Item(13)
$ws.Range("D7").Value = 0.01414105531172678
$ws.Range("E7").Value = 0.1005107813526911
$ws.Range("F7").Value = 0.1522337718945985
$ws.Range("G7").Value = 0.1189161692610672
$ws.Range("H7").Value = 15.98104812649535
$ws.Range("L7").Value = 0.01326207254649109
$ws.Range("M7").Value = 0.09794732433530426
$ws.Range("N7").Value = 0.1461621916531723
$ws.Range("O7").Value = 0.1151610721836641
$ws.Range("P7").Value = 15.12229778952996

$ws = $wb.Worksheets.Item(14)
$ws.Range("D7").Value = 0.01117562372515061
$ws.Range("E7").Value = 0.08775522437540648
$ws.Range("F7").Value = 0.1635377920128878
$ws.Range("G7").Value = 0.1057148226368971
$ws.Range("H7").Value = 14.7779894258248
$ws.Range("L7").Value = 0.01292078822757069
$ws.Range("M7").Value = 0.09405278915623863
$ws.Range("N7").Value = 0.1728796094902846
$ws.Range("O7").Value = 0.1136696451457938
$ws.Range("P7").Value = 16.95508660090787

$ws = $wb.Worksheets.Item(15)
$ws.Range("D7").Value = 0.02221693215870389
$ws.Range("E7").Value = 0.1131962566343289
$ws.Range("F7").Value = 0.286437217716897
$ws.Range("G7").Value = 0.1490534540314443
$ws.Range("H7").Value = 20.63423305811549
$ws.Range("L7").Value = 0.02086729007015207
$ws.Range("M7").Value = 0.1070461573852048
$ws.Range("N7").Value = 0.2385449401851052
$ws.Range("O7").Value = 0.1444551489914848
$ws.Range("P7").Value = 19.21986675165265

$ws = $wb.Worksheets.Item(16)
$ws.Range("D7").Value = 0.1424954309215581
$ws.Range("E7").Value = 0.3238790400138878
$ws.Range("F7").Value = 0.1194523839454943
$ws.Range("G7").Value = 0.3774856698227869
$ws.Range("H7").Value = 12.52955403119356
$ws.Range("L7").Value = 0.1533346021161947
$ws.Range("M7").Value = 0.3362770004627986
$ws.Range("N7").Value = 0.1233110045643023
$ws.Range("O7").Value = 0.3915796242352182
$ws.Range("P7").Value = 12.97615244636943

$ws = $wb.Worksheets.Item(17)
$ws.Range("D7").Value = 0.1370007494831095
$ws.Range("E7").Value = 0.2670387841134421
$ws.Range("F7").Value = 0.1044484430432899
$ws.Range("G7").Value = 0.3701361229103551
$ws.Range("H7").Value = 10.79012910228141
$ws.Range("L7").Value = 0.1153728697050223
$ws.Range("M7").Value = 0.2463381412108982
$ws.Range("N7").Value = 0.09989899465728802
$ws.Range("O7").Value = 0.3396658206311349
$ws.Range("P7").Value = 9.974543686142995

$ws = $wb.Worksheets.Item(18)
$ws.Range("D7").Value = 0.1427271061563565
$ws.Range("E7").Value = 0.263515967714446
$ws.Range("F7").Value = 0.1368934471464912
$ws.Range("G7").Value = 0.3777924114594634
$ws.Range("H7").Value = 11.79432360949121
$ws.Range("L7").Value = 0.1717288377725348
$ws.Range("M7").Value = 0.3118884371422596
$ws.Range("N7").Value = 0.1434178081059787
$ws.Range("O7").Value = 0.4144017830228711
$ws.Range("P7").Value = 13.38303864170244

$ws = $wb.Worksheets.Item(19)
$ws.Range("D7").Value = 0.07553508556925363
$ws.Range("E7").Value = 0.233315061967965
$ws.Range("F7").Value = 0.3598025590105745
$ws.Range("G7").Value = 0.2748364705952498
$ws.Range("H7").Value = 33.304122457757

$ws = $wb.Worksheets.Item(20)
$ws.Range("D7").Value = 0.3468894770984437
$ws.Range("E7").Value = 0.453215162776795
$ws.Range("F7").Value = 0.9776808523268423
$ws.Range("G7").Value = 0.5889732397133538
$ws.Range("H7").Value = 64.10376008316366

$ws = $wb.Worksheets.Item(21)
$ws.Range("D7").Value = 0.08609101259860948
$ws.Range("E7").Value = 0.2316553556493845
$ws.Range("F7").Value = 0.4818173560327427
$ws.Range("G7").Value = 0.2934127001317589
$ws.Range("H7").Value = 36.12027165921437

$ws = $wb.Worksheets.Item(22)
$ws.Range("D7").Value = 4.158700297774248
$ws.Range("E7").Value = 1.508370200784822
$ws.Range("F7").Value = 0.5591964328158379
$ws.Range("G7").Value = 2.039289164825393
$ws.Range("H7").Value = 74.39153803359714

$ws = $wb.Worksheets.Item(23)
$ws.Range("D7").Value = 14.05013415788708
$ws.Range("E7").Value = 3.455321587307366
$ws.Range("F7").Value = 1.291517959153331
$ws.Range("G7").Value = 3.748350858429222
$ws.Range("H7").Value = 130.6492280357874

$ws = $wb.Worksheets.Item(24)
$ws.Range("D7").Value = 10.6770172545855
$ws.Range("E7").Value = 2.764809982734926
$ws.Range("F7").Value = 1.298046434133867
$ws.Range("G7").Value = 3.267570543168961
$ws.Range("H7").Value = 110.2914527976326
